$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New log entry in row 14: task description, date, and hours "geboden" (offered)
$ws.Range("A14").Value = "klant kan user aan duiden op taak, taak waar een toewijzing op gebeurd kan niet meer geboden op worden en komen achteraan op de bord te staan. "
$ws.Range("B14").Value = (Get-Date -Year 2024 -Month 12 -Day 15).Date
$ws.Range("D14").Value = 1.5

# Row grows to fit the wrapped text (4 lines @ 15.75pt)
$ws.Rows(14).RowHeight = 63

# Move the view: scrolled to row 13, active cell on the new entry
$win = $excel.ActiveWindow
$win.ScrollRow = 13
$win.ScrollColumn = 1
$ws.Range("A14").Select()
